# Insert a new data row at row 214 (pushes the existing row 214..307 block
# down to 215..308, matching the target dimension A1:R308) and populate the
# new row with the inserted record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 214:307 down to 215:308, leaving a fresh blank row 214.
$ws.Rows.Item(214).Insert()

# Populate the new row 214 with the new record (Berenjena @ Terminal La
# Palmera de La Serena, Coquimbo).
$ws.Range("A214").Value = 8
$ws.Range("B214").Value = "Terminal La Palmera de La Serena"
$ws.Range("C214").Value = "Coquimbo"
$ws.Range("D214").Value = 45205
$ws.Range("E214").Value = 4
$ws.Range("F214").Value = 100112001
$ws.Range("G214").Value = "Berenjena"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 400
$ws.Range("K214").Value = 9500
$ws.Range("L214").Value = 10000
$ws.Range("M214").Value = 9750
$ws.Range("N214").Value = "`$/caja 50 unidades"
$ws.Range("O214").Value = "Región de Arica y Parinacota"
$ws.Range("P214").Value = 195
$ws.Range("Q214").Value = 50
$ws.Range("R214").Value = "Hortaliza"
